$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24 (current "Music/mp3"), shifting rows 24:33 down to 25:34.
$ws.Rows("24").Insert()

# Populate the newly inserted row 24 with "Videos" / "gif".
$ws.Range("A24").Value = "Videos"
$ws.Range("B24").Value = "gif"

# Grow the ConfigTable (ListObject) to include the newly added row.
$table = $ws.ListObjects.Item("ConfigTable")
$table.Resize($ws.Range("A1:B34"))
